$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting rows 3-11 down to 4-12
$ws.Rows.Item(3).Insert()

# Fill in the new row 3 with the "spear" spell data
$ws.Range("A3").Value = "spear"
$ws.Range("B3").Value = "attack"
$ws.Range("C3").Value = "a basic long cool atk"
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 90
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = "MT"

# Update selection to H3
$ws.Range("H3").Select()
